$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 208.5
$ws.Range("I2").Value = 174.9
$ws.Range("J2").Value = 292.5
$ws.Range("K2").Value = 174.9
$ws.Range("L2").Value = 292.5
$ws.Range("M2").Value = -61.90000000000001
$ws.Range("N2").Value = -518.5
$ws.Range("H4").Value = 136.75
$ws.Range("I4").Value = 136.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 136.75
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -22.75
$ws.Range("H38").Value = 565
$ws.Range("I38").Value = 14.166667
$ws.Range("J38").Value = 1666.6666
$ws.Range("K38").Value = 42.500001
$ws.Range("L38").Value = 4999.9998
$ws.Range("M38").Value = 329.499999
$ws.Range("N38").Value = -5743.9998
$ws.Range("H39").Value = 244.6875
$ws.Range("I39").Value = 69.75
$ws.Range("J39").Value = 419.625
$ws.Range("K39").Value = 209.25
$ws.Range("L39").Value = 1258.875
$ws.Range("M39").Value = 86.75
$ws.Range("N39").Value = -1850.875
$ws.Range("H40").Value = 2450.6428
$ws.Range("I40").Value = 2407.1428
$ws.Range("J40").Value = 2494.1428
$ws.Range("K40").Value = 2407.1428
$ws.Range("L40").Value = 2494.1428
$ws.Range("M40").Value = -2232.1428
$ws.Range("N40").Value = -2844.1428
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("M46").Value = -29881
$ws.Range("N46").Value = -30238
$ws.Range("H60").Value = 10000
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 10000
$ws.Range("K60").Value = 30000
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -29516
$ws.Range("N60").Value = -30968
$ws.Range("H138").Value = 6862198
$ws.Range("I138").Value = 1840405.4
$ws.Range("J138").Value = 10419301
$ws.Range("K138").Value = 5521216.199999999
$ws.Range("L138").Value = 31257903
$ws.Range("M138").Value = -5516076.199999999
$ws.Range("N138").Value = -31268183
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 47499.8
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 47499.8
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 47499.8
$ws.Range("N133").Value = -52559.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 211158.33
$ws.Range("I105").Value = 2579.0881
$ws.Range("J105").Value = 717707.9399999999
$ws.Range("K105").Value = 2579.0881
$ws.Range("L105").Value = 717707.9399999999
$ws.Range("M105").Value = -832.0880999999999
$ws.Range("N105").Value = -721201.9399999999
$ws.Range("H139").Value = 90780
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 90780
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 90780
$ws.Range("N139").Value = -101060
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 40000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 40000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40496
$ws.Range("H67").Value = 40000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 40000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41716
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 7000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 7000
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").Value = 21000
$ws.Range("N102").Value = -25868
$ws.Range("H122").Value = 1077.04
$ws.Range("I122").Value = 375.33334
$ws.Range("J122").Value = 1471.75
$ws.Range("K122").Value = 3378.00006
$ws.Range("L122").Value = 13245.75
$ws.Range("M122").Value = -928.0000600000003
$ws.Range("N122").Value = -18145.75
$ws.Range("H123").Value = 2507
$ws.Range("I123").Value = 1515
$ws.Range("J123").Value = 3499
$ws.Range("K123").Value = 4545
$ws.Range("L123").Value = 10497
$ws.Range("M123").Value = -2095
$ws.Range("N123").Value = -15397
$ws.Range("H131").Value = 1302.1625
$ws.Range("I131").Value = 424.2857
$ws.Range("J131").Value = 1386.3424
$ws.Range("K131").Value = 1272.8571
$ws.Range("L131").Value = 4159.0272
$ws.Range("M131").Value = 3767.1429
$ws.Range("N131").Value = -14239.0272
$ws.Range("H132").Value = 1463.2
$ws.Range("I132").Value = 1334
$ws.Range("J132").Value = 1576.25
$ws.Range("K132").Value = 12006
$ws.Range("L132").Value = 14186.25
$ws.Range("M132").Value = -9476
$ws.Range("N132").Value = -19246.25
$ws.Range("H141").Value = 3922.3076
$ws.Range("I141").Value = 3699.0908
$ws.Range("J141").Value = 5150
$ws.Range("K141").Value = 11097.2724
$ws.Range("L141").Value = 15450
$ws.Range("M141").Value = -5917.2724
$ws.Range("N141").Value = -25810
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16003.429
$ws.Range("I43").Value = 16256
$ws.Range("J43").Value = 15666.667
$ws.Range("K43").Value = 16256
$ws.Range("L43").Value = 15666.667
$ws.Range("M43").Value = -16105
$ws.Range("N43").Value = -15968.667
$ws.Range("H80").Value = 3286.6667
$ws.Range("I80").Value = 3063.6365
$ws.Range("J80").Value = 3900
$ws.Range("K80").Value = 3063.6365
$ws.Range("L80").Value = 3900
$ws.Range("M80").Value = -2065.6365
$ws.Range("N80").Value = -5896
$ws.Range("H83").Value = 3286.6667
$ws.Range("I83").Value = 3063.6365
$ws.Range("J83").Value = 3900
$ws.Range("K83").Value = 15318.1825
$ws.Range("L83").Value = 19500
$ws.Range("M83").Value = -10326.1825
$ws.Range("N83").Value = -29484
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 415.88235
$ws.Range("I55").Value = 296.83334
$ws.Range("J55").Value = 701.6
$ws.Range("K55").Value = 296.83334
$ws.Range("L55").Value = 701.6
$ws.Range("M55").Value = -123.83334
$ws.Range("N55").Value = -1047.6
$ws.Range("H82").Value = 1508.8636
$ws.Range("I82").Value = 1345.9375
$ws.Range("J82").Value = 1943.3334
$ws.Range("K82").Value = 1345.9375
$ws.Range("L82").Value = 1943.3334
$ws.Range("M82").Value = -984.9375
$ws.Range("N82").Value = -2665.3334
$ws.Range("H85").Value = 1508.8636
$ws.Range("I85").Value = 1345.9375
$ws.Range("J85").Value = 1943.3334
$ws.Range("K85").Value = 1345.9375
$ws.Range("L85").Value = 1943.3334
$ws.Range("M85").Value = -97.9375
$ws.Range("N85").Value = -4439.3334
$ws.Range("H132").Value = 3848.9355
$ws.Range("I132").Value = 2179.85
$ws.Range("J132").Value = 6883.636
$ws.Range("K132").Value = 6539.549999999999
$ws.Range("L132").Value = 20650.908
$ws.Range("M132").Value = -4009.549999999999
$ws.Range("N132").Value = -25710.908
$ws.Range("H133").Value = 57624.57
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 57624.57
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 57624.57
$ws.Range("N133").Value = -62684.57
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2007.2812
$ws.Range("I132").Value = 1413.6666
$ws.Range("J132").Value = 3788.125
$ws.Range("K132").Value = 4240.9998
$ws.Range("L132").Value = 11364.375
$ws.Range("M132").Value = -1710.9998
$ws.Range("N132").Value = -16424.375
$ws.Range("H135").Value = 294444
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 294444
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 294444
$ws.Range("N135").Value = -304584
